# "Add files via upload" – the workbook had been authored in a PT-locale
# Excel, so its only worksheet still carried the default name "Folha1"
# ("Sheet1"). Re-uploading it under its real name means the tab should be
# renamed to match the workbook's subject: "DummyData".
$wb = $excel.ActiveWorkbook

$sheet = $wb.Sheets.Item(1)
$sheet.Name = "DummyData"
